# Proyecto de Clase. Ajuste BBDD
# Updates the USUARIOS, MENSAJES and CREDENCIALES sample-data tables on the
# MER worksheet, and moves the current viewport/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S.I.W.Ventas")

# --- USUARIOS table: row 21 (customer-2) changes owner from Pepito Perez to Jorge Campos
$ws.Range("E21").Value = "Jorge"
$ws.Range("F21").Value = "Campos"

# --- MENSAJES table: row 27 (was person-2 / "Reunión selleres")
$ws.Range("C27").Value = "seller-1"
$ws.Range("E27").Value = "Reunión Vendedores"
$ws.Range("F27").Value = "El próximo fin de semana…"

# --- MENSAJES table: row 28 description text
$ws.Range("F28").Value = "Qué papeles piden para…"

# --- MENSAJES table: row 30 (date + asunto + descripcion)
$ws.Range("D30").Value = 44805
$ws.Range("E30").Value = "No tengo el producto"
$ws.Range("F30").Value = "Aun no llega el producto…"

# --- MENSAJES table: row 29 (was customer-1 / "No tengo el producto")
$ws.Range("C29").Value = "admin-1"
$ws.Range("E29").Value = "Mantenimiento Sistema"
$ws.Range("F29").Value = "Se informa a los usuarios…"

# --- CREDENCIALES table: row 38 (was admin-2)
$ws.Range("C38").Value = "customer-2"

# --- Viewport / selection change
$ws.Application.ActiveWindow.ScrollRow = 24
$ws.Range("A24").Select()
$ws.Range("D38").Select()
